$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a "text" value into a cell using the existing text-number-format
# style (cellXfs index 1 in the original workbook: numFmtId=49 applyNumberFormat=1),
# without introducing any new styles. We force the NumberFormat to Text ("@") so the
# digit-only strings are NOT coerced into numeric cells, then reset back to the
# workbook's "Normal" cell style (which collapses back onto the existing numFmt=49
# style because the font/fill/border are unchanged) and finally re-apply "@" so the
# xf ends up identical to the pre-existing style index 1.
function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.Style = "Normal"
    $r.NumberFormat = "@"
}

# Row 3: O3:T3 new numeric truth-table cells
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 1
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 1
$ws.Range("T3").Value = 0

# Row 4: O4:T4
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = 0
$ws.Range("Q4").Value = 1
$ws.Range("R4").Value = 0
$ws.Range("S4").Value = 1
$ws.Range("T4").Value = 1

# Row 5: O5:T5
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = 0
$ws.Range("Q5").Value = 1
$ws.Range("R5").Value = 0
$ws.Range("S5").Value = 0
$ws.Range("T5").Value = 1

# Row 6: O6:T6
$ws.Range("O6").Value = 0
$ws.Range("P6").Value = 0
$ws.Range("Q6").Value = 1
$ws.Range("R6").Value = 1
$ws.Range("S6").Value = 0
$ws.Range("T6").Value = 0

# Row 7: O7:T7
$ws.Range("O7").Value = 0
$ws.Range("P7").Value = 0
$ws.Range("Q7").Value = 1
$ws.Range("R7").Value = 1
$ws.Range("S7").Value = 0
$ws.Range("T7").Value = 1

# Row 8: O8:T8
$ws.Range("O8").Value = 0
$ws.Range("P8").Value = 0
$ws.Range("Q8").Value = 1
$ws.Range("R8").Value = 1
$ws.Range("S8").Value = 1
$ws.Range("T8").Value = 0

# Row 12: empty text-formatted cell
$ws.Range("O12").NumberFormat = "@"

# Row 13: text value "101011" (opcode column)
Set-TextValue "O13" "101011"

# Row 14: empty text-formatted cell
$ws.Range("O14").NumberFormat = "@"

# Row 15: H15 = "0010", O15 = "001011"
Set-TextValue "H15" "0010"
Set-TextValue "O15" "001011"

# Row 16: G16 empty text-formatted cell, H16 = "1000"
$ws.Range("G16").NumberFormat = "@"
Set-TextValue "H16" "1000"

# Row 17: H17 = "1010"
Set-TextValue "H17" "1010"

# Row 18: H18 = "1100"
Set-TextValue "H18" "1100"

# Row 2: new bold blank cell Q2 (new font+style added here)
$ws.Range("Q2").Font.Bold = $true

# Row 22 / Row 23: move the "1" from F -> E (keeping the existing style on E),
# then drop F entirely.
$ws.Range("E22").Value = 1
$ws.Range("F22").Clear()

$ws.Range("E23").Value = 1
$ws.Range("F23").Clear()

# Update the active selection to match the saved view state.
$ws.Range("L15").Select()
